$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns keep their original text format,
# since many of the new values look like plain numbers (e.g. "1.001", "317.11")
# and Excel would otherwise silently convert them to numeric cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.195.35'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.806.22'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.91%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.39%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.11'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.25%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5335'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3782'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.45%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07482'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.03'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.11%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.214'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.73%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.369'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.808.55'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.86%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.74'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001064'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.16%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.22%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.39'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.88%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.37%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.226.69'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.27%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.085'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.07'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.51'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.014.55'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.323'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.14'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.113'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1088'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +7.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.582'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.70%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.622'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07197'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +9.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2229'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02298'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.458'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6182'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.88%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.442'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.180'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.45'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.687'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5769'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.38'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.45%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.186'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.925'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.49%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06823'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.85'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.05%  '
